# "deleted TS loves me"
# Remove the paragraph containing "I LOVE TS SHE LOVES ME" (run text +
# its bookmarkStart/bookmarkEnd) entirely, including its paragraph mark,
# so the following (empty) paragraph takes its place.

$d = $word.ActiveDocument

# Find the paragraph whose text contains the target line and delete the
# whole paragraph (text + paragraph mark) by expanding the range up to
# the start of the next paragraph.
$paras = $d.Paragraphs
for ($i = $paras.Count; $i -ge 1; $i--) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*I LOVE TS SHE LOVES ME*") {
        $nextStart = $d.Content.End
        if ($i -lt $paras.Count) {
            $nextStart = $paras.Item($i + 1).Range.Start
        }
        $r = $d.Range($p.Range.Start, $nextStart)
        $r.Delete()
    }
}
